# NATMI LR-pair output refresh: "update scripts wuth new tpm"
#
# The per-cell-type TPM table feeding this Jag2 -> Notch4 ligand-receptor
# sheet was regenerated, which changes the per-cluster ligand/receptor
# average & total expression values and every value derived from them
# (derived specificities + edge weights) for all 16 sending x target
# cluster combinations. Columns A:F and K:L (identities, cell counts,
# detection rates) are unaffected by the new TPM and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("G2").Value = 20.81798233333333
$ws.Range("H2").Value = 62.453947
$ws.Range("I2").Value = 0.8242653639952813
$ws.Range("J2").Value = 0.8242653639952813
$ws.Range("M2").Value = 46.63275166666667
$ws.Range("N2").Value = 139.898255
$ws.Range("O2").Value = 0.9158911059585902
$ws.Range("P2").Value = 0.9158911059585902
$ws.Range("Q2").Value = 970.7998003513873
$ws.Range("R2").Value = 8737.198203162485
$ws.Range("S2").Value = 0.7549373158329982
$ws.Range("T2").Value = 0.7549373158329982

# Row 3: ECs -> FAPs
$ws.Range("G3").Value = 20.81798233333333
$ws.Range("H3").Value = 62.453947
$ws.Range("I3").Value = 0.8242653639952813
$ws.Range("J3").Value = 0.8242653639952813
$ws.Range("O3").Value = 0.05441917700612491
$ws.Range("P3").Value = 0.05441917700612491
$ws.Range("Q3").Value = 57.68166742654389
$ws.Range("R3").Value = 519.1350068388949
$ws.Range("S3").Value = 0.04485584274327718
$ws.Range("T3").Value = 0.04485584274327719

# Row 4: ECs -> MuSCs
$ws.Range("G4").Value = 20.81798233333333
$ws.Range("H4").Value = 62.453947
$ws.Range("I4").Value = 0.8242653639952813
$ws.Range("J4").Value = 0.8242653639952813
$ws.Range("M4").Value = 0.849605
$ws.Range("N4").Value = 2.548815
$ws.Range("O4").Value = 0.01668667696558362
$ws.Range("P4").Value = 0.01668667696558362
$ws.Range("Q4").Value = 17.68706188031167
$ws.Range("R4").Value = 159.183556922805
$ws.Range("S4").Value = 0.01375424986290846
$ws.Range("T4").Value = 0.01375424986290846

# Row 5: ECs -> Resolving-Mac
$ws.Range("G5").Value = 20.81798233333333
$ws.Range("H5").Value = 62.453947
$ws.Range("I5").Value = 0.8242653639952813
$ws.Range("J5").Value = 0.8242653639952813
$ws.Range("M5").Value = 0.662052
$ws.Range("N5").Value = 1.986156
$ws.Range("O5").Value = 0.01300304006970129
$ws.Range("P5").Value = 0.0130030400697013
$ws.Range("Q5").Value = 13.782586839748
$ws.Range("R5").Value = 124.043281557732
$ws.Range("S5").Value = 0.01071795555609757
$ws.Range("T5").Value = 0.01071795555609757

# Row 6: FAPs -> ECs
$ws.Range("G6").Value = 0.7925996666666667
$ws.Range("I6").Value = 0.03138212158540782
$ws.Range("J6").Value = 0.03138212158540782
$ws.Range("M6").Value = 46.63275166666667
$ws.Range("N6").Value = 139.898255
$ws.Range("O6").Value = 0.9158911059585902
$ws.Range("P6").Value = 0.9158911059585902
$ws.Range("Q6").Value = 36.96110342674945
$ws.Range("R6").Value = 332.649930840745
$ws.Range("S6").Value = 0.02874260604618611
$ws.Range("T6").Value = 0.02874260604618611

# Row 7: FAPs -> FAPs
$ws.Range("G7").Value = 0.7925996666666667
$ws.Range("I7").Value = 0.03138212158540782
$ws.Range("J7").Value = 0.03138212158540782
$ws.Range("O7").Value = 0.05441917700612491
$ws.Range("P7").Value = 0.05441917700612491
$ws.Range("Q7").Value = 2.196104773412778
$ws.Range("R7").Value = 19.764942960715
$ws.Range("S7").Value = 0.001707789229384041
$ws.Range("T7").Value = 0.001707789229384041

# Row 8: FAPs -> MuSCs
$ws.Range("G8").Value = 0.7925996666666667
$ws.Range("I8").Value = 0.03138212158540782
$ws.Range("J8").Value = 0.03138212158540782
$ws.Range("M8").Value = 0.849605
$ws.Range("N8").Value = 2.548815
$ws.Range("O8").Value = 0.01668667696558362
$ws.Range("P8").Value = 0.01668667696558362
$ws.Range("Q8").Value = 0.6733966397983334
$ws.Range("R8").Value = 6.060569758185
$ws.Range("S8").Value = 0.0005236633253903693
$ws.Range("T8").Value = 0.0005236633253903693

# Row 9: FAPs -> Resolving-Mac
$ws.Range("G9").Value = 0.7925996666666667
$ws.Range("I9").Value = 0.03138212158540782
$ws.Range("J9").Value = 0.03138212158540782
$ws.Range("M9").Value = 0.662052
$ws.Range("N9").Value = 1.986156
$ws.Range("O9").Value = 0.01300304006970129
$ws.Range("P9").Value = 0.0130030400697013
$ws.Range("Q9").Value = 0.524742194516
$ws.Range("R9").Value = 4.722679750644
$ws.Range("S9").Value = 0.0004080629844472958
$ws.Range("T9").Value = 0.0004080629844472958

# Row 10: MuSCs -> ECs
$ws.Range("G10").Value = 1.536855
$ws.Range("H10").Value = 4.610564999999999
$ws.Range("I10").Value = 0.06085010188305478
$ws.Range("J10").Value = 0.06085010188305479
$ws.Range("M10").Value = 46.63275166666667
$ws.Range("N10").Value = 139.898255
$ws.Range("O10").Value = 0.9158911059585902
$ws.Range("P10").Value = 0.9158911059585902
$ws.Range("Q10").Value = 71.667777562675
$ws.Range("R10").Value = 645.0099980640749
$ws.Range("S10").Value = 0.05573206711136394
$ws.Range("T10").Value = 0.05573206711136394

# Row 11: MuSCs -> FAPs
$ws.Range("G11").Value = 1.536855
$ws.Range("H11").Value = 4.610564999999999
$ws.Range("I11").Value = 0.06085010188305478
$ws.Range("J11").Value = 0.06085010188305479
$ws.Range("O11").Value = 0.05441917700612491
$ws.Range("P11").Value = 0.05441917700612491
$ws.Range("Q11").Value = 4.258258921224999
$ws.Range("R11").Value = 38.32433029102499
$ws.Range("S11").Value = 0.003311412465214692
$ws.Range("T11").Value = 0.003311412465214693

# Row 12: MuSCs -> MuSCs
$ws.Range("G12").Value = 1.536855
$ws.Range("H12").Value = 4.610564999999999
$ws.Range("I12").Value = 0.06085010188305478
$ws.Range("J12").Value = 0.06085010188305479
$ws.Range("M12").Value = 0.849605
$ws.Range("N12").Value = 2.548815
$ws.Range("O12").Value = 0.01668667696558362
$ws.Range("P12").Value = 0.01668667696558362
$ws.Range("Q12").Value = 1.305719692275
$ws.Range("R12").Value = 11.751477230475
$ws.Range("S12").Value = 0.001015385993445387
$ws.Range("T12").Value = 0.001015385993445387

# Row 13: MuSCs -> Resolving-Mac
$ws.Range("G13").Value = 1.536855
$ws.Range("H13").Value = 4.610564999999999
$ws.Range("I13").Value = 0.06085010188305478
$ws.Range("J13").Value = 0.06085010188305479
$ws.Range("M13").Value = 0.662052
$ws.Range("N13").Value = 1.986156
$ws.Range("O13").Value = 0.01300304006970129
$ws.Range("P13").Value = 0.0130030400697013
$ws.Range("Q13").Value = 1.01747792646
$ws.Range("R13").Value = 9.157301338139998
$ws.Range("S13").Value = 0.0007912363130307674
$ws.Range("T13").Value = 0.0007912363130307676

# Row 14: Resolving-Mac -> ECs
$ws.Range("G14").Value = 2.108971
$ws.Range("H14").Value = 6.326912999999999
$ws.Range("I14").Value = 0.08350241253625613
$ws.Range("J14").Value = 0.08350241253625615
$ws.Range("M14").Value = 46.63275166666667
$ws.Range("N14").Value = 139.898255
$ws.Range("O14").Value = 0.9158911059585902
$ws.Range("P14").Value = 0.9158911059585902
$ws.Range("Q14").Value = 98.34712091520167
$ws.Range("R14").Value = 885.1240882368149
$ws.Range("S14").Value = 0.07647911696804208
$ws.Range("T14").Value = 0.07647911696804209

# Row 15: Resolving-Mac -> FAPs
$ws.Range("G15").Value = 2.108971
$ws.Range("H15").Value = 6.326912999999999
$ws.Range("I15").Value = 0.08350241253625613
$ws.Range("J15").Value = 0.08350241253625615
$ws.Range("O15").Value = 0.05441917700612491
$ws.Range("P15").Value = 0.05441917700612491
$ws.Range("Q15").Value = 5.843456002911666
$ws.Range("R15").Value = 52.59110402620499
$ws.Range("S15").Value = 0.004544132568248986
$ws.Range("T15").Value = 0.004544132568248987

# Row 16: Resolving-Mac -> MuSCs
$ws.Range("G16").Value = 2.108971
$ws.Range("H16").Value = 6.326912999999999
$ws.Range("I16").Value = 0.08350241253625613
$ws.Range("J16").Value = 0.08350241253625615
$ws.Range("M16").Value = 0.849605
$ws.Range("N16").Value = 2.548815
$ws.Range("O16").Value = 0.01668667696558362
$ws.Range("P16").Value = 0.01668667696558362
$ws.Range("Q16").Value = 1.791792306455
$ws.Range("R16").Value = 16.126130758095
$ws.Range("S16").Value = 0.001393377783839406
$ws.Range("T16").Value = 0.001393377783839407

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Range("G17").Value = 2.108971
$ws.Range("H17").Value = 6.326912999999999
$ws.Range("I17").Value = 0.08350241253625613
$ws.Range("J17").Value = 0.08350241253625615
$ws.Range("M17").Value = 0.662052
$ws.Range("N17").Value = 1.986156
$ws.Range("O17").Value = 0.01300304006970129
$ws.Range("P17").Value = 0.0130030400697013
$ws.Range("Q17").Value = 1.396248468492
$ws.Range("R17").Value = 12.566236216428
$ws.Range("S17").Value = 0.001085785216125666
$ws.Range("T17").Value = 0.001085785216125667

